$d = $word.ActiveDocument

# 1. Merge the three runs in "User Stories #0" story paragraph into a single
#    run by re-finding/replacing the full sentence (Word collapses the
#    matched run fragments that share identical formatting into one run).
$d.Content.Find.Execute(
    "Eu Carl Miller enquanto gestor de infraestrutura, desejo ser capaz de identificar que colaborador está alocado com qual máquina além de poder alterar esses registros para assim ter um melhor manejo sobre os patrimônios da empresa;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Eu Carl Miller enquanto gestor de infraestrutura, desejo ser capaz de identificar que colaborador está alocado com qual máquina além de poder alterar esses registros para assim ter um melhor manejo sobre os patrimônios da empresa;",
    2) | Out-Null

# 2. Merge the three runs in the "User Stories #5" story paragraph into a
#    single run the same way.
$d.Content.Find.Execute(
    "Eu Carl Miller enquanto gestor de infraestrutura, gostaria que o sistema além de monitorar as máquinas, também realizasse tarefas rotineiras para manter as máquinas em bom estado;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Eu Carl Miller enquanto gestor de infraestrutura, gostaria que o sistema além de monitorar as máquinas, também realizasse tarefas rotineiras para manter as máquinas em bom estado;",
    2) | Out-Null

# 3. Remove the 11 trailing empty paragraphs that followed that last story
#    paragraph, so the section properties immediately follow it again.
$trailingStart = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "também realizasse tarefas rotineiras") {
        $trailingStart = $i + 1
        break
    }
}

if ($trailingStart -ne $null -and $trailingStart -le $d.Paragraphs.Count) {
    $startRange = $d.Paragraphs.Item($trailingStart).Range.Start
    $endRange = $d.Content.End
    $r = $d.Range($startRange, $endRange)
    $r.Delete()
}
